## Test results now use format testNumer_+testNumber.ToString+workflowName+now
## Adds a third data row (row 3) to the Tests sheet: same WorkflowFile as row 2,
## with a run-count of 7, widens column A to fit the longer workflow path, and
## moves the active selection onto the newly written cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of test results: reuse the existing workflow file name already in A2
# and record 7 runs for it in B3.
$ws.Range("A3").Value = "Test_Framework\Tests\_wbTest_Example1.xaml"
$ws.Range("B3").Value = 7

# Column A needs to be a bit wider now to comfortably show the workflow path;
# column B keeps its existing width, so leave it untouched.
$ws.Columns.Item(1).ColumnWidth = 57.7109375

# Leave the selection on the cell we just filled in.
[void]$ws.Range("B3").Select()
